# Append a new "Action items" row ("Assessment Criterias" / "" /
# "Aim for Wednesday") to the end of the Action-items table, right
# after the existing "Presentation Plan" row.

$d = $word.ActiveDocument

# Locate the target table robustly: the one whose last row's first
# cell contains "Presentation Plan" (rather than hard-coding a table
# index).
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $lastRow = $t.Rows.Item($t.Rows.Count)
    if ($lastRow.Cells.Item(1).Range.Text -like "*Presentation Plan*") {
        $targetTable = $t
    }
}

if ($targetTable -eq $null) {
    throw "Could not locate the Action items table (Presentation Plan row not found)."
}

# Add a new row at the end of the table and fill in its cells.
$newRow = $targetTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Assessment Criterias"
$newRow.Cells.Item(3).Range.Text = "Aim for Wednesday"
